# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# sheets, which carry duplicate data for the same events.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F14" = 1540
    "F16" = 489
    "F22" = 1409
    "F28" = 78
    "F29" = 1776
    "F36" = 638
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
